# Updates cryptos list values per upstream data refresh (GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text formatted like "1.234.56"; some new values
# look like plain decimals (e.g. "309.29") which Excel would otherwise auto-convert
# to a number. Force those specific cells to Text format first so they stay strings,
# matching the original inline-string cell type.
$ws.Range("D2").Value = '44.556.50'
$ws.Range("E2").Value = '  +3.49%  '
$ws.Range("D3").Value = '2.431.39'
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.29'
$ws.Range("E5").Value = '  +2.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.98'
$ws.Range("E6").Value = '  +5.21%  '
$ws.Range("E7").Value = '  +1.57%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.38'
$ws.Range("E10").Value = '  +3.40%  '
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("E13").Value = '  +2.63%  '
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("D15").Value = '2.814.53'
$ws.Range("E15").Value = '  +2.38%  '
$ws.Range("D16").Value = '2.418.69'
$ws.Range("E16").Value = '  +1.71%  '
$ws.Range("E17").Value = '  +4.35%  '
$ws.Range("D18").Value = '44.457.96'
$ws.Range("E18").Value = '  +3.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.46'
$ws.Range("E19").Value = '  +2.40%  '
$ws.Range("E20").Value = '  +1.52%  '
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.87'
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.33'
$ws.Range("E23").Value = '  +4.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.59'
$ws.Range("E24").Value = '  +2.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("E25").Value = '  +2.59%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.18'
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("E28").Value = '  -2.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.68'
$ws.Range("E29").Value = '  +4.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.19'
$ws.Range("E30").Value = '  +5.35%  '
$ws.Range("E31").Value = '  +16.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.63'
$ws.Range("E32").Value = '  +12.35%  '
$ws.Range("E33").Value = '  +2.40%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("E35").Value = '  +3.29%  '
$ws.Range("E36").Value = '  +1.95%  '
$ws.Range("E37").Value = '  +3.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.89'
$ws.Range("E38").Value = '  +3.04%  '
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '126.23'
$ws.Range("E40").Value = '  +11.64%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.109'
$ws.Range("E41").Value = '  +0.78%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.89'
$ws.Range("E42").Value = '  -2.68%  '
$ws.Range("E43").Value = '  +3.70%  '
$ws.Range("D44").Value = '1.947.99'
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("E46").Value = '  +6.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.56'
$ws.Range("E47").Value = '  +4.05%  '
$ws.Range("E48").Value = '  +10.26%  '
$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.47'
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.97'
$ws.Range("E50").Value = '  +2.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.63'
$ws.Range("E51").Value = '  +4.58%  '
